# Added filtering options for the Component Analysis
#
# The forecast-error table on Sheet1 previously had every horizon column
# (G..K, i.e. Q5..Q9) populated for every origin row. The evaluation
# window actually available for each origin quarter is limited, so the
# "future" horizon cells that fall outside of the available evaluation
# window need to be cleared out (emptied), leaving the staircase pattern
# of populated cells that is already visible further down the sheet
# (rows 45 and below).
#
# Clear exactly the trailing cells identified in the target diff, row by
# row, using a single multi-area Range and looping over its Areas so
# every discontiguous block gets cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rangesToClear = "G2:K2,I3:K3,G4:K4,I5:K5,G6:K6,I7:K7,G8:K8,I9:K9,G10:K10,I11:K11,G12:K12,I13:K13,G14:K14,I15:K15,G16:K16,I17:K17,K18,I19:K19,K20,I21:K21,K22,J23:K23,I24:K24,K26,J27:K27,I28:K28,K30,J31:K31,I32:K32,K34,J35:K35,I36:K36,K38,J39:K39,I40:K40,K42,J43:K43,I44:J44"

$clearRange = $ws.Range($rangesToClear)
foreach ($area in $clearRange.Areas) {
    $area.ClearContents()
}
